# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets, which contain identical data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1501
    $ws.Range("F6").Value = 40
    $ws.Range("F9").Value = 302
}
